$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8301
$ws1.Range("F3").Value = 7714
$ws1.Range("F4").Value = 116
$ws1.Range("F9").Value = 112
$ws1.Range("F10").Value = 157
$ws1.Range("F12").Value = 698
$ws1.Range("F14").Value = 1276
$ws1.Range("F16").Value = 48
$ws1.Range("F19").Value = 107

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8301
$ws4.Range("F3").Value = 7714
$ws4.Range("F4").Value = 116
$ws4.Range("F9").Value = 112
$ws4.Range("F10").Value = 157
$ws4.Range("F12").Value = 698
$ws4.Range("F14").Value = 1276
$ws4.Range("F16").Value = 48
$ws4.Range("F19").Value = 107
